# fix createdAt mockup data
# Adds "createdAt" / "updatedAt" columns (L, M) to the user mockup sheet,
# populating each data row (2-20) with a timestamp (same value used for
# both createdAt and updatedAt, matching the source mockup data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers
$ws.Range("L1").Value = "createdAt"
$ws.Range("M1").Value = "updatedAt"

# One timestamp per data row (row 2 .. row 20)
$timestamps = @(
  "2023-04-16T03:42:46.836Z",
  "2023-04-15T21:32:21.365Z",
  "2023-04-17T15:18:00.971Z",
  "2023-04-20T18:59:51.241Z",
  "2023-04-12T20:57:30.058Z",
  "2023-04-18T13:51:56.194Z",
  "2023-04-18T14:12:52.457Z",
  "2023-04-20T16:38:07.175Z",
  "2023-04-25T03:04:07.411Z",
  "2023-04-21T19:02:23.154Z",
  "2023-04-27T03:56:47.284Z",
  "2023-05-01T06:19:23.798Z",
  "2023-05-01T00:47:26.804Z",
  "2023-04-29T21:45:55.351Z",
  "2023-05-01T15:47:35.774Z",
  "2023-05-06T15:09:27.627Z",
  "2023-05-06T20:36:21.009Z",
  "2023-05-15T10:22:09.736Z",
  "2023-05-10T16:58:08.913Z"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 12).Value = $timestamps[$i]
  $ws.Cells.Item($row, 13).Value = $timestamps[$i]
}

# Move the active selection (matches the author's final cursor position)
$ws.Range("I7").Select()
